$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (L1:N1), matching the existing header style ---
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- Fix E/F columns (percent-like values were scaled by 100) ---
$eVals = @(87.71535580524345, 12.28464419475655, 85.55858310626702, 14.44141689373297, 95.90643274853801, 4.093567251461988)
$fVals = @(61.31511528608027, 68.29268292682927, 93.31210191082803, 99.52830188679245, 21.79878048780488, 32.14285714285715)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Range("E$row").Value = $eVals[$i]
    $ws.Range("F$row").Value = $fVals[$i]
}

# --- New data columns L, M, N for rows 2-7 ---
$lVals = @(92.54049761504434, 81.11912441944573, 91.96147942061606, 78.3393122619878, 19.45553018286073, 24.5476882057119)
$mVals = @(235006, 28547, 176818, 26828, 2163, 45)
$nVals = @(327.3064066852368, 254.8839285714286, 150.8686006825939, 127.1469194312796, 15.12587412587413, 5)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Range("L$row").Value = $lVals[$i]
    $ws.Range("M$row").Value = $mVals[$i]
    $ws.Range("N$row").Value = $nVals[$i]
}
